# statistics.xlsx edit: add the 2016-08-25 ("day 18") data row, let it
# flow through the existing rolling-average / error formulas (same
# pattern every prior row uses), and leave the workbook with the Chart1
# sheet active (matching the author's final view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- new raw data for row 18 (2016-08-25) --------------------------------
$ws.Cells.Item(18, 2).Value = 7     # B18 总访问数
$ws.Cells.Item(18, 3).Value = 0     # C18 Kalles Fraktaler

# --- prediction columns E:K (same pattern as every prior row) -----------
# E=avg(prev 1), F=avg(prev 2), G=avg(prev 3), ... K=avg(prev 7)
$predFormulas = @{
    5  = '=(B17+B18)/2'
    6  = '=(B16+B18)/2'
    7  = '=(B15+B18)/2'
    8  = '=(B14+B18)/2'
    9  = '=(B13+B18)/2'
    10 = '=(B12+B18)/2'
    11 = '=(B11+B18)/2'
}
foreach ($col in $predFormulas.Keys) {
    $c = $ws.Cells.Item(18, $col)
    $c.Formula = $predFormulas[$col]
    $c.NumberFormat = "0.0"
    $c.VerticalAlignment = -4108   # xlVAlignCenter
}

# --- spacer column L (formatted, no value) -------------------------------
$l18 = $ws.Cells.Item(18, 12)
$l18.NumberFormat = "0.00"
$l18.VerticalAlignment = -4108

# --- error columns M:S (abs difference vs each prediction) --------------
$errFormulas = @{
    13 = '=ABS($B18-E18)'
    14 = '=ABS($B18-F18)'
    15 = '=ABS($B18-G18)'
    16 = '=ABS($B18-H18)'
    17 = '=ABS($B18-I18)'
    18 = '=ABS($B18-J18)'
    19 = '=ABS($B18-K18)'
}
foreach ($col in $errFormulas.Keys) {
    $c = $ws.Cells.Item(18, $col)
    $c.Formula = $errFormulas[$col]
    $c.NumberFormat = "0.0"
    $c.VerticalAlignment = -4108
}

# --- leave Chart1 as the active/selected sheet, like the saved file -----
$chartSheet = $wb.Worksheets.Item("Chart1")
$chartSheet.Activate()
